# Auto-generated edit script: updates horarios-141 scrape data across 3 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 10:32:07"
$ws.Cells.Item(3,1).Value = "Total filas: 111"
$ws.Cells.Item(15,3).Value = "215A_EL PATO"
$ws.Cells.Item(16,3).Value = "225_GOMEZ"
$ws.Cells.Item(23,1).Value = "06:46:40"
$ws.Cells.Item(23,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(23,4).Value = 35
$ws.Cells.Item(24,1).Value = "06:15:23"
$ws.Cells.Item(24,3).Value = "16_SANTA ANA"
$ws.Cells.Item(24,4).Value = 66
$ws.Cells.Item(35,1).Value = "07:26:49"
$ws.Cells.Item(35,3).Value = "16_SANTA ANA"
$ws.Cells.Item(35,4).Value = 34
$ws.Cells.Item(36,1).Value = "07:51:40"
$ws.Cells.Item(36,3).Value = "17_ROMERO"
$ws.Cells.Item(36,4).Value = 9
$ws.Cells.Item(43,1).Value = "08:14:55"
$ws.Cells.Item(43,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(43,4).Value = 14
$ws.Cells.Item(44,1).Value = "06:58:58"
$ws.Cells.Item(44,3).Value = "15_ABASTO"
$ws.Cells.Item(44,4).Value = 90
$ws.Cells.Item(45,1).Value = "07:51:40"
$ws.Cells.Item(45,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(45,4).Value = 38
$ws.Cells.Item(46,1).Value = "08:14:55"
$ws.Cells.Item(46,3).Value = "15_ABASTO"
$ws.Cells.Item(46,4).Value = 15
$ws.Cells.Item(55,1).Value = "08:49:06"
$ws.Cells.Item(55,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(55,4).Value = 4
$ws.Cells.Item(56,1).Value = "08:14:55"
$ws.Cells.Item(56,3).Value = "215B_EL PATO"
$ws.Cells.Item(56,4).Value = 39
$ws.Cells.Item(71,1).Value = "08:14:55"
$ws.Cells.Item(71,3).Value = "16_SANTA ANA"
$ws.Cells.Item(71,4).Value = 77
$ws.Cells.Item(72,1).Value = "08:49:06"
$ws.Cells.Item(72,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(72,4).Value = 42
$ws.Cells.Item(90,3).Value = "14_ABASTO"
$ws.Cells.Item(91,3).Value = "15_ABASTO"
$ws.Cells.Item(92,1).Value = "10:32:07"
$ws.Cells.Item(92,4).Value = 12
$ws.Cells.Item(93,1).Value = "10:32:07"
$ws.Cells.Item(93,4).Value = 14
$ws.Cells.Item(95,1).Value = "10:32:07"
$ws.Cells.Item(95,2).Value = "10:55"
$ws.Cells.Item(95,3).Value = "16_SANTA ANA"
$ws.Cells.Item(95,4).Value = 23
$ws.Cells.Item(96,1).Value = "10:32:07"
$ws.Cells.Item(96,2).Value = "10:57"
$ws.Cells.Item(96,3).Value = "27_EL RETIRO"
$ws.Cells.Item(96,4).Value = 25
$ws.Cells.Item(97,1).Value = "10:32:07"
$ws.Cells.Item(97,2).Value = "10:59"
$ws.Cells.Item(97,3).Value = "10_OLMOS"
$ws.Cells.Item(97,4).Value = 27
$ws.Cells.Item(98,1).Value = "10:32:07"
$ws.Cells.Item(98,2).Value = "11:01"
$ws.Cells.Item(98,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(98,4).Value = 29
$ws.Cells.Item(99,1).Value = "10:32:07"
$ws.Cells.Item(99,2).Value = "11:06"
$ws.Cells.Item(99,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(99,4).Value = 34
$ws.Cells.Item(100,1).Value = "10:32:07"
$ws.Cells.Item(100,2).Value = "11:10"
$ws.Cells.Item(100,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(100,4).Value = 38
$ws.Cells.Item(101,2).Value = "11:14"
$ws.Cells.Item(101,3).Value = "14_ABASTO"
$ws.Cells.Item(101,4).Value = 92
$ws.Cells.Item(102,1).Value = "10:32:07"
$ws.Cells.Item(102,2).Value = "11:15"
$ws.Cells.Item(102,3).Value = "15X38_ABASTO"
$ws.Cells.Item(102,4).Value = 43
$ws.Cells.Item(103,1).Value = "10:32:07"
$ws.Cells.Item(103,2).Value = "11:15"
$ws.Cells.Item(103,3).Value = "14_ABASTO"
$ws.Cells.Item(103,4).Value = 43
$ws.Cells.Item(103,5).Value = "LP1912"
$ws.Cells.Item(104,1).Value = "10:32:07"
$ws.Cells.Item(104,2).Value = "11:24"
$ws.Cells.Item(104,3).Value = "16_SANTA ANA"
$ws.Cells.Item(104,4).Value = 52
$ws.Cells.Item(104,5).Value = "LP1912"
$ws.Cells.Item(105,1).Value = "10:32:07"
$ws.Cells.Item(105,2).Value = "11:29"
$ws.Cells.Item(105,3).Value = "10_OLMOS"
$ws.Cells.Item(105,4).Value = 57
$ws.Cells.Item(105,5).Value = "LP1912"
$ws.Cells.Item(106,1).Value = "09:42:42"
$ws.Cells.Item(106,2).Value = "11:30"
$ws.Cells.Item(106,3).Value = "215C_EL PATO"
$ws.Cells.Item(106,4).Value = 108
$ws.Cells.Item(106,5).Value = "LP1912"
$ws.Cells.Item(107,1).Value = "10:32:07"
$ws.Cells.Item(107,2).Value = "11:31"
$ws.Cells.Item(107,3).Value = "215C_EL PATO"
$ws.Cells.Item(107,4).Value = 59
$ws.Cells.Item(107,5).Value = "LP1912"
$ws.Cells.Item(108,1).Value = "10:32:07"
$ws.Cells.Item(108,2).Value = "11:42"
$ws.Cells.Item(108,3).Value = "215B_EL PATO"
$ws.Cells.Item(108,4).Value = 70
$ws.Cells.Item(108,5).Value = "LP1912"
$ws.Cells.Item(109,1).Value = "10:32:07"
$ws.Cells.Item(109,2).Value = "11:45"
$ws.Cells.Item(109,3).Value = "15X38_ABASTO"
$ws.Cells.Item(109,4).Value = 73
$ws.Cells.Item(109,5).Value = "LP1912"
$ws.Cells.Item(110,1).Value = "10:32:07"
$ws.Cells.Item(110,2).Value = "11:51"
$ws.Cells.Item(110,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(110,4).Value = 79
$ws.Cells.Item(110,5).Value = "LP1912"
$ws.Cells.Item(111,1).Value = "10:32:07"
$ws.Cells.Item(111,2).Value = "11:53"
$ws.Cells.Item(111,3).Value = "225_GOMEZ"
$ws.Cells.Item(111,4).Value = 81
$ws.Cells.Item(111,5).Value = "LP1912"
$ws.Cells.Item(112,1).Value = "10:32:07"
$ws.Cells.Item(112,2).Value = "11:58"
$ws.Cells.Item(112,3).Value = "17_ROMERO"
$ws.Cells.Item(112,4).Value = 86
$ws.Cells.Item(112,5).Value = "LP1912"
$ws.Cells.Item(113,1).Value = "10:32:07"
$ws.Cells.Item(113,2).Value = "12:06"
$ws.Cells.Item(113,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(113,4).Value = 94
$ws.Cells.Item(113,5).Value = "LP1912"
$ws.Cells.Item(114,1).Value = "10:32:07"
$ws.Cells.Item(114,2).Value = "12:10"
$ws.Cells.Item(114,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(114,4).Value = 98
$ws.Cells.Item(114,5).Value = "LP1912"
$ws.Cells.Item(115,1).Value = "10:32:07"
$ws.Cells.Item(115,2).Value = "12:10"
$ws.Cells.Item(115,3).Value = "15_ABASTO"
$ws.Cells.Item(115,4).Value = 98
$ws.Cells.Item(115,5).Value = "LP1912"
$ws.Cells.Item(116,1).Value = "10:32:07"
$ws.Cells.Item(116,2).Value = "12:22"
$ws.Cells.Item(116,3).Value = "215C_EL PATO"
$ws.Cells.Item(116,4).Value = 110
$ws.Cells.Item(116,5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 10:32:07"
$ws.Cells.Item(3,1).Value = "Total filas: 20"
$ws.Cells.Item(23,1).Value = "10:32:07"
$ws.Cells.Item(23,2).Value = "11:31"
$ws.Cells.Item(23,3).Value = "215C_EL PATO"
$ws.Cells.Item(23,4).Value = 59
$ws.Cells.Item(23,5).Value = "LP1912"
$ws.Cells.Item(24,1).Value = "10:32:07"
$ws.Cells.Item(24,2).Value = "11:42"
$ws.Cells.Item(24,3).Value = "215B_EL PATO"
$ws.Cells.Item(24,4).Value = 70
$ws.Cells.Item(24,5).Value = "LP1912"
$ws.Cells.Item(25,1).Value = "10:32:07"
$ws.Cells.Item(25,2).Value = "12:22"
$ws.Cells.Item(25,3).Value = "215C_EL PATO"
$ws.Cells.Item(25,4).Value = 110
$ws.Cells.Item(25,5).Value = "LP1912"

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 10:32:07"
$ws.Cells.Item(3,1).Value = "Total filas: 20"
$ws.Cells.Item(19,1).Value = "09:42:42"
$ws.Cells.Item(19,3).Value = "215A_LA PLATA"
$ws.Cells.Item(19,4).Value = 48
$ws.Cells.Item(20,1).Value = "08:49:06"
$ws.Cells.Item(20,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(20,4).Value = 101
$ws.Cells.Item(22,1).Value = "10:32:07"
$ws.Cells.Item(22,2).Value = "10:32"
$ws.Cells.Item(22,3).Value = "215A_LA PLATA"
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = "L6173"
$ws.Cells.Item(23,1).Value = "10:32:07"
$ws.Cells.Item(23,2).Value = "10:32"
$ws.Cells.Item(23,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = "L6173"
$ws.Cells.Item(24,1).Value = "09:42:42"
$ws.Cells.Item(24,2).Value = "11:25"
$ws.Cells.Item(24,3).Value = "215C_LA PLATA"
$ws.Cells.Item(24,4).Value = 103
$ws.Cells.Item(24,5).Value = "L6203"
$ws.Cells.Item(25,1).Value = "10:32:07"
$ws.Cells.Item(25,2).Value = "11:26"
$ws.Cells.Item(25,3).Value = "215C_LA PLATA"
$ws.Cells.Item(25,4).Value = 54
$ws.Cells.Item(25,5).Value = "L6203"
